$wb = $excel.ActiveWorkbook

# --- Build the new TODOS/COMBINADAS data (28 rows, cols A:D) ---
$dataTodos = New-Object 'object[,]' 28,4
$dataTodos[0,0] = "16:48"
$dataTodos[0,1] = "16_SANTA ANA"
$dataTodos[0,2] = 0
$dataTodos[0,3] = "🚌"
$dataTodos[1,0] = "16:48"
$dataTodos[1,1] = "14_ABASTO"
$dataTodos[1,2] = 0
$dataTodos[1,3] = "🚌"
$dataTodos[2,0] = "16:48"
$dataTodos[2,1] = "16_SANTA ANA"
$dataTodos[2,2] = 0
$dataTodos[2,3] = "🚌"
$dataTodos[3,0] = "16:49"
$dataTodos[3,1] = "215B_LP-P MOR-40 Y 115"
$dataTodos[3,2] = 1
$dataTodos[3,3] = "🚌"
$dataTodos[4,0] = "16:55"
$dataTodos[4,1] = "23_HERNANDEZ"
$dataTodos[4,2] = 7
$dataTodos[4,3] = "🚌"
$dataTodos[5,0] = "16:56"
$dataTodos[5,1] = "17_179 Y 38"
$dataTodos[5,2] = 8
$dataTodos[5,3] = "🚌"
$dataTodos[6,0] = "16:57"
$dataTodos[6,1] = "10_OLMOS"
$dataTodos[6,2] = 9
$dataTodos[6,3] = "🚌"
$dataTodos[7,0] = "17:04"
$dataTodos[7,1] = "11_ETCHEVERRY"
$dataTodos[7,2] = 16
$dataTodos[7,3] = "🚌"
$dataTodos[8,0] = "17:04"
$dataTodos[8,1] = "215A_EL PATO"
$dataTodos[8,2] = 16
$dataTodos[8,3] = "🚌"
$dataTodos[9,0] = "17:04"
$dataTodos[9,1] = "23_HERNANDEZ"
$dataTodos[9,2] = 16
$dataTodos[9,3] = "🚌"
$dataTodos[10,0] = "17:10"
$dataTodos[10,1] = "10_OLMOS"
$dataTodos[10,2] = 22
$dataTodos[10,3] = "🚌"
$dataTodos[11,0] = "17:14"
$dataTodos[11,1] = "215A_LA PLATA"
$dataTodos[11,2] = 26
$dataTodos[11,3] = "🚌"
$dataTodos[12,0] = "17:16"
$dataTodos[12,1] = "11_ETCHEVERRY"
$dataTodos[12,2] = 28
$dataTodos[12,3] = "🚌"
$dataTodos[13,0] = "17:21"
$dataTodos[13,1] = "26_HERNANDEZ"
$dataTodos[13,2] = 33
$dataTodos[13,3] = "🚌"
$dataTodos[14,0] = "17:28"
$dataTodos[14,1] = "14_ABASTO"
$dataTodos[14,2] = 40
$dataTodos[14,3] = "🚌"
$dataTodos[15,0] = "17:36"
$dataTodos[15,1] = "15_ABASTO"
$dataTodos[15,2] = 48
$dataTodos[15,3] = "🚌"
$dataTodos[16,0] = "17:36"
$dataTodos[16,1] = "27_EL RETIRO"
$dataTodos[16,2] = 48
$dataTodos[16,3] = "🚌"
$dataTodos[17,0] = "17:38"
$dataTodos[17,1] = "17_ROMERO"
$dataTodos[17,2] = 50
$dataTodos[17,3] = "📅"
$dataTodos[18,0] = "17:40"
$dataTodos[18,1] = "215B_EL PATO"
$dataTodos[18,2] = 52
$dataTodos[18,3] = "📅"
$dataTodos[19,0] = "17:45"
$dataTodos[19,1] = "15_ABASTO"
$dataTodos[19,2] = 57
$dataTodos[19,3] = "🚌"
$dataTodos[20,0] = "17:50"
$dataTodos[20,1] = "16_P MOR-167 Y 521"
$dataTodos[20,2] = 62
$dataTodos[20,3] = "🚌"
$dataTodos[21,0] = "17:52"
$dataTodos[21,1] = "81_EL PELIGRO"
$dataTodos[21,2] = 64
$dataTodos[21,3] = "📅"
$dataTodos[22,0] = "18:03"
$dataTodos[22,1] = "215C_LA PLATA"
$dataTodos[22,2] = 75
$dataTodos[22,3] = "🚌"
$dataTodos[23,0] = "18:04"
$dataTodos[23,1] = "17_ROMERO"
$dataTodos[23,2] = 76
$dataTodos[23,3] = "🚌"
$dataTodos[24,0] = "18:14"
$dataTodos[24,1] = "84_COLONIA URQUIZA-ESC 49"
$dataTodos[24,2] = 86
$dataTodos[24,3] = "🚌"
$dataTodos[25,0] = "18:21"
$dataTodos[25,1] = "26_HERNANDEZ"
$dataTodos[25,2] = 93
$dataTodos[25,3] = "🚌"
$dataTodos[26,0] = "18:27"
$dataTodos[26,1] = "215C_EL PATO"
$dataTodos[26,2] = 99
$dataTodos[26,3] = "🚌"
$dataTodos[27,0] = "18:32"
$dataTodos[27,1] = "11X44_ETCHEVERRY"
$dataTodos[27,2] = 104
$dataTodos[27,3] = "🚌"

# --- Build the new 215 data (6 rows, cols A:D) ---
$data215 = New-Object 'object[,]' 6,4
$data215[0,0] = "16:49"
$data215[0,1] = "215B_LP-P MOR-40 Y 115"
$data215[0,2] = 1
$data215[0,3] = "🚌"
$data215[1,0] = "17:04"
$data215[1,1] = "215A_EL PATO"
$data215[1,2] = 16
$data215[1,3] = "🚌"
$data215[2,0] = "17:14"
$data215[2,1] = "215A_LA PLATA"
$data215[2,2] = 26
$data215[2,3] = "🚌"
$data215[3,0] = "17:40"
$data215[3,1] = "215B_EL PATO"
$data215[3,2] = 52
$data215[3,3] = "📅"
$data215[4,0] = "18:03"
$data215[4,1] = "215C_LA PLATA"
$data215[4,2] = 75
$data215[4,3] = "🚌"
$data215[5,0] = "18:27"
$data215[5,1] = "215C_EL PATO"
$data215[5,2] = 99
$data215[5,3] = "🚌"

# --- Apply to TODOS sheet ---
$wsTodos = $wb.Worksheets.Item("TODOS")
$wsTodos.Range("A2:D29").Value = $dataTodos
$wsTodos.Rows.Item(31).Delete()
$wsTodos.Rows.Item(30).Delete()

# --- Apply to 215 sheet (row count unchanged: A1:D7) ---
$ws215 = $wb.Worksheets.Item("215")
$ws215.Range("A2:D7").Value = $data215

# --- Apply to COMBINADAS sheet (same shape/content as TODOS) ---
$wsComb = $wb.Worksheets.Item("COMBINADAS")
$wsComb.Range("A2:D29").Value = $dataTodos
$wsComb.Rows.Item(31).Delete()
$wsComb.Rows.Item(30).Delete()

